$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Potrebno ispraviti" table (rows 2-9) ---
# Row 6 ("Dijagram kod opisa sustava ..."): assign to Katarina, mark as done (DA)
$ws.Range("C6").Value = "Katarina"
$ws.Range("D6").Value = "DA"

# Row 7 ("Arhitektura: opis razreda ..."): note that the fix is only partial
$ws.Range("F7").Value = "DJELOMIČNO: dodani uvod i opis modela i contollera"

# --- "2. dio dokumentacije" table (rows 12-22) ---
# Row 19 ("Ispitivanje programskog rješenja"): mark as not done (NE)
$ws.Range("D19").Value = "NE"

# Row 21 ("Korisničke upute"): assign to Ante
$ws.Range("C21").Value = "Ante"

# Update the sheet view's selection to C15 (top-left scroll resets to default)
$ws.Activate()
$ws.Range("C15").Select()
